# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from serial date 45185 (2023-09-16) to serial date 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
